$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "638×5=" "834×3="
Replace-Text "338×9=" "248×9="
Replace-Text "380×9=" "689×8="
Replace-Text "693×9=" "131×8="
Replace-Text "914×4=" "785×7="
Replace-Text "629×8=" "561×8="
Replace-Text "976×9=" "194×4="
Replace-Text "117×7=" "706×3="
Replace-Text "387×6=" "751×4="
Replace-Text "805×8=" "548×7="
Replace-Text "771×4=" "956×6="
Replace-Text "136×7=" "968×2="
Replace-Text "507×7=" "680×4="
Replace-Text "110×2=" "807×7="
Replace-Text "730×9=" "509×4="
Replace-Text "241×8=" "222×5="
Replace-Text "248×7=" "803×6="
Replace-Text "820×4=" "318×3="
Replace-Text "835×8=" "312×8="
Replace-Text "384×2=" "653×2="
Replace-Text "958×8=" "357×6="
Replace-Text "261×2=" "838×3="
Replace-Text "329×7=" "544×2="
Replace-Text "523×5=" "954×6="
Replace-Text "855×3=" "679×4="
